$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "IsPriority" column (L) to the server-config table, mirroring the
# formatting of column K (the preceding column) for each of the header /
# type / data rows, then fill in the new column's content.
$ws.Range("K3:K7").Copy($ws.Range("L3:L7"))

$ws.Range("L3").Value = "是否默认值"
$ws.Range("L4").Value = "IsPriority"
$ws.Range("L5").Value = "int"
$ws.Range("L6").Value = 1
$ws.Range("L7").Value = 0

# Reflect the new active cell/selection left behind after the edit.
$ws.Range("K12").Select()
